$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "CasesTab" Cypher query in B2: the `Cohort` field/match was
#    removed from the query (it referenced an undefined/erroneous cohort
#    traversal), per the commit "Fixed variables and query errors".
$b2Text = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Saint Bernard']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@
$ws.Range("B2").Value = $b2Text

# 2. Row heights shrink now that the B2 query text is shorter (fewer
#    wrapped lines); also rows 3/4 get re-measured slightly.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 216
$ws.Rows.Item(4).RowHeight = 216

# 3. Reset the view: scroll back to the top-left (was parked at A4) and
#    move the selection to B2.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
